$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '94.491.29'
$ws.Range('E2').Value = '  -3.33%  '
$ws.Range('D3').Value = '3.434.41'
$ws.Range('E3').Value = '  +2.40%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = '''237.90'
$ws.Range('E5').Value = '  -5.37%  '
$ws.Range('D6').Value = '''642.33'
$ws.Range('E6').Value = '  -2.40%  '
$ws.Range('D7').Value = '''1.45'
$ws.Range('E7').Value = '  +4.02%  '
$ws.Range('D8').Value = '''0.405'
$ws.Range('E8').Value = '  -4.05%  '
$ws.Range('E9').Value = '  +0.14%  '
$ws.Range('D10').Value = '''0.973'
$ws.Range('E10').Value = '  -2.85%  '
$ws.Range('D11').Value = '3.435.12'
$ws.Range('E11').Value = '  +2.47%  '
$ws.Range('D12').Value = '''41.94'
$ws.Range('E12').Value = '  +2.23%  '
$ws.Range('D13').Value = '''0.198'
$ws.Range('E13').Value = '  -4.82%  '
$ws.Range('D14').Value = '''6.20'
$ws.Range('E14').Value = '  +1.99%  '
$ws.Range('D15').Value = '94.342.81'
$ws.Range('E15').Value = '  -3.11%  '
$ws.Range('D16').Value = '4.076.48'
$ws.Range('E16').Value = '  +2.35%  '
$ws.Range('D17').Value = '''0.0000251'
$ws.Range('E17').Value = '  -1.15%  '
$ws.Range('D18').Value = '''8.26'
$ws.Range('E18').Value = '  -4.78%  '
$ws.Range('D19').Value = '3.434.25'
$ws.Range('E19').Value = '  +2.59%  '
$ws.Range('D20').Value = '''17.58'
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('D21').Value = '''11.40'
$ws.Range('E21').Value = '  +5.22%  '
$ws.Range('D22').Value = '''0.510'
$ws.Range('E22').Value = '  -1.91%  '
$ws.Range('D23').Value = '''499.04'
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('D24').Value = '''3.22'
$ws.Range('E24').Value = '  -3.71%  '
$ws.Range('D25').Value = '''0.0000192'
$ws.Range('E25').Value = '  -3.74%  '
$ws.Range('D26').Value = '''6.52'
$ws.Range('E26').Value = '  -4.81%  '
$ws.Range('D27').Value = '''94.52'
$ws.Range('E27').Value = '  +1.27%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '3.619.41'
$ws.Range('E28').Value = '  +2.69%  '
$ws.Range('B29').Value = 'Aptos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D29').Value = '''11.91'
$ws.Range('E29').Value = '  -2.47%  '
$ws.Range('D30').Value = '''11.68'
$ws.Range('E30').Value = '  +3.11%  '
$ws.Range('B31').Value = 'Dai'
$ws.Range('C31').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D31').Value = '''1.00'
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = '''2.77'
$ws.Range('E32').Value = '  +7.92%  '
$ws.Range('D33').Value = '''0.138'
$ws.Range('E33').Value = '  -2.11%  '
$ws.Range('E34').Value = '  -0.51%  '
$ws.Range('D35').Value = '''0.177'
$ws.Range('E35').Value = '  -4.60%  '
$ws.Range('D36').Value = '''29.85'
$ws.Range('E36').Value = '  +5.02%  '
$ws.Range('D37').Value = '''0.553'
$ws.Range('E37').Value = '  -0.69%  '
$ws.Range('D38').Value = '''548.58'
$ws.Range('E38').Value = '  +4.11%  '
$ws.Range('D39').Value = '''7.64'
$ws.Range('E39').Value = '  -4.25%  '
$ws.Range('D40').Value = '''1.45'
$ws.Range('E40').Value = '  -4.46%  '
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '''0.151'
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').Value = '''0.900'
$ws.Range('E43').Value = '  +6.75%  '
$ws.Range('D44').Value = '''24.04'
$ws.Range('E44').Value = '  -1.52%  '
$ws.Range('D45').Value = '''1.71'
$ws.Range('E45').Value = '  -1.38%  '
$ws.Range('D46').Value = '''5.57'
$ws.Range('E46').Value = '  +0.21%  '
$ws.Range('B47').Value = 'dogwifhat'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D47').Value = '''3.33'
$ws.Range('E47').Value = '  +4.98%  '
$ws.Range('D48').Value = '''0.0409'
$ws.Range('E48').Value = '  -3.81%  '
$ws.Range('B49').Value = 'MantraDAO'
$ws.Range('C49').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D49').Value = '''3.57'
$ws.Range('E49').Value = '  -0.29%  '
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').Value = '''54.23'
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = '''2.16'
$ws.Range('E51').Value = '  -5.49%  '

Write-Output "done"
